$wb = $excel.ActiveWorkbook

# --- "Platform Coverage" sheet: fill in the yearly coverage row ---
$wsCov = $wb.Worksheets.Item("Platform Coverage")

# Row 2 already had 0.6 in every other year-column from H (2018) through
# AD (2040). Fill in the remaining in-between year columns with the same
# 0.6 value so every year from H2 to AD2 has a coverage value defined.
for ($col = 8; $col -le 30; $col++) {
    $wsCov.Cells.Item(2, $col).Value = 0.6
}

# --- Make "Platform Coverage" the active sheet / selection ---
$wsCov.Activate()

# Scroll the view so column Q is the left-most visible column (best
# effort — mirrors the intent of the saved view state).
try {
    $excel.ActiveWindow.ScrollColumn = 17
} catch {
}

$wsCov.Range("AC2").Select()

# --- "MarketShare" sheet: leave its own last selection (Y6) untouched.
#     It is no longer the active/selected tab now that Platform Coverage
#     has been activated above. ---
